# Atualizei dados bibi e add
# Update row 9 (Ano 2025) figures in the faturamento_anual worksheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = 3702697.74
$ws.Range("C9").Value = 582980.24
$ws.Range("D9").Value = 4285677.98
$ws.Range("E9").Value = 13.60298750210812
$ws.Range("F9").Value = 86.39701249789186
$ws.Range("G9").Value = -43.65786667748834
$ws.Range("H9").Value = -33.13439683109686
$ws.Range("I9").Value = 37363
$ws.Range("J9").Value = 1590
$ws.Range("K9").Value = 38953
$ws.Range("L9").Value = 26914
$ws.Range("M9").Value = 159.2360102548859
$ws.Range("N9").Value = 8.713667476979792
